# [client] more on german translation
#
# Cleans up three placeholder lines in the "Geschätzter Wert der
# Liegenschaft / Wohnrecht / Abschlagzahlung" block: each line used to end
# with a trailing space plus a run of tab characters (and, for the first
# line, two extra trailing spaces) intended to leave room for a
# hand-written value. The edit trims all of that so each line ends
# cleanly right after the colon.

$d = $word.ActiveDocument

function Trim-TrailingTabs($searchText, $replaceText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, `
                                $false, $true, 1, $false, $replaceText, 2)
    if ($found) {
        # The Find/Replace above rewrites the run's text nodes and, as a
        # side effect, drops the (empty) <w:rPr/> element that every run
        # in this document otherwise carries. Touching a character
        # formatting property on the freshly-replaced range (and putting
        # it straight back) makes the run keep an explicit - if empty -
        # run-properties element, matching the rest of the document.
        $rng2 = $d.Range($rng.Start, $rng.Start + $replaceText.Length)
        $rng2.Font.Bold = $true
        $rng2.Font.Bold = $false
    }
    return $found
}

# "Geschätzter Wert der Liegenschaft: " + three tabs + two spaces -> clean colon
Trim-TrailingTabs "Geschätzter Wert der Liegenschaft: `t`t`t  " "Geschätzter Wert der Liegenschaft:" | Out-Null

# "Wohnrecht: " + six tabs -> clean colon
Trim-TrailingTabs "Wohnrecht: `t`t`t`t`t`t" "Wohnrecht:" | Out-Null

# "Abschlagzahlung: " + five tabs -> clean colon
Trim-TrailingTabs "Abschlagzahlung: `t`t`t`t`t" "Abschlagzahlung:" | Out-Null
